$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells that are no longer used in the new layout ---
# Old "new"/"new damage" block header cells that moved down to row 8/9
$ws.Range("Q1").ClearContents()
$ws.Range("A3:A7").ClearContents()
$ws.Range("T2:U7").ClearContents()
$ws.Range("X2:Y7").ClearContents()
$ws.Range("Z1:AD7").ClearContents()
$ws.Range("AF1:AG7").ClearContents()

# --- Write the new layout ---
$ws.Range("C1").Value = "old"
$ws.Range("L1").Value = "old damage"
$ws.Range("V1").Value = "irl stats"
$ws.Range("A2").Value = "name"
$ws.Range("B2").Value = "pretty_name"
$ws.Range("C2").Value = "ergonomicsold"
$ws.Range("D2").Value = "weightold"
$ws.Range("E2").Value = "horizontal_recoilold"
$ws.Range("F2").Value = "vertical_recoilold"
$ws.Range("G2").Value = "bullet_deviationold"
$ws.Range("H2").Value = "bullet_damageold"
$ws.Range("I2").Value = "bullet_velocityold"
$ws.Range("J2").Value = "fire_rateold"
$ws.Range("K2").Value = "priceold"
$ws.Range("L2").Value = "0st"
$ws.Range("M2").Value = "100st"
$ws.Range("N2").Value = "200st"
$ws.Range("O2").Value = "avg"
$ws.Range("P2").Value = "vel loss"
$ws.Range("Q2").Value = "suppression"
$ws.Range("R2").Value = "pen"
$ws.Range("S2").Value = "strength"
$ws.Range("V2").Value = "mv"
$ws.Range("W2").Value = "energy"
$ws.Range("B3").Value = "5.56x45 Federal American Eagle Training 55gr FMJ (XM193BLX)"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.1
$ws.Range("G3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 52
$ws.Range("M3").Value = 48
$ws.Range("N3").Value = 42
$ws.Range("O3").Formula = "=AVERAGE(L3:N3)"
$ws.Range("P3").Value = 40
$ws.Range("Q3").Value = 20
$ws.Range("R3").Value = 0.91
$ws.Range("S3").Formula = "=C3-D3*20-E3*0.8-F3-0.6-G3*5+I3/200+(O3-50)*1.5"
$ws.Range("V3").Value = 3165
$ws.Range("W3").Value = 1658.165
$ws.Range("B4").Value = "5.56x45 CBC Defense M196 55gr Red Tracer"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0.1
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.05
$ws.Range("I4").Value = 10
$ws.Range("K4").Value = 600
$ws.Range("L4").Value = 52
$ws.Range("M4").Value = 48
$ws.Range("N4").Value = 42
$ws.Range("O4").Formula = "=AVERAGE(L4:N4)"
$ws.Range("P4").Value = 40
$ws.Range("Q4").Value = 30
$ws.Range("R4").Value = 0.8
$ws.Range("S4").Formula = "=C4-D4*20-E4*0.8-F4-0.6-G4*5+I4/200+(O4-50)*1.5"
$ws.Range("V4").Value = 3200
$ws.Range("W4").Value = 1512
$ws.Range("B5").Value = "5.56x45 Winchester M855 62gr FMJ Green Tip (VM855)"
$ws.Range("C5").Value = -1
$ws.Range("D5").Value = 0.13
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = -0.1
$ws.Range("I5").Value = 50
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 58
$ws.Range("M5").Value = 55
$ws.Range("N5").Value = 51
$ws.Range("O5").Formula = "=AVERAGE(L5:N5)"
$ws.Range("P5").Value = 50
$ws.Range("Q5").Value = 100
$ws.Range("R5").Value = 0.63
$ws.Range("S5").Formula = "=C5-D5*20-E5*0.8-F5-0.6-G5*5+I5/200+(O5-50)*1.5"
$ws.Range("V5").Value = 3060
$ws.Range("W5").Value = 1748
$ws.Range("B6").Value = "5.56x45 CBC Defense SAT IP 62gr CLF (AEP-97)"
$ws.Range("C6").Value = -3
$ws.Range("D6").Value = 0.17
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = -0.4
$ws.Range("I6").Value = 300
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 52
$ws.Range("M6").Value = 48
$ws.Range("N6").Value = 42
$ws.Range("O6").Formula = "=AVERAGE(L6:N6)"
$ws.Range("P6").Value = 20
$ws.Range("Q6").Value = 10
$ws.Range("R6").Value = 0.97
$ws.Range("S6").Formula = "=C6-D6*20-E6*0.8-F6-0.6-G6*5+I6/200+(O6-50)*1.5"
$ws.Range("V6").Value = 2953
$ws.Range("W6").Value = 1620
$ws.Range("B7").Value = "5.56x45 NOVX Copper Pentagon 55gr MCHP (556N55CP)"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0.07
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 0.3
$ws.Range("I7").Value = -100
$ws.Range("K7").Value = 1200
$ws.Range("L7").Value = 62
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 38
$ws.Range("O7").Formula = "=AVERAGE(L7:N7)"
$ws.Range("P7").Value = 80
$ws.Range("Q7").Value = 80
$ws.Range("R7").Value = 0.38
$ws.Range("S7").Formula = "=C7-D7*20-E7*0.8-F7-0.6-G7*5+I7/200+(O7-50)*1.5"
$ws.Range("V7").Value = 3340
$ws.Range("W7").Value = 1844
$ws.Range("C8").Value = "new"
$ws.Range("L8").Value = "new damage"
$ws.Range("C9").Value = "ergonomics"
$ws.Range("D9").Value = "weight"
$ws.Range("E9").Value = "horizontal_recoil"
$ws.Range("F9").Value = "vertical_recoil"
$ws.Range("G9").Value = "bullet_deviation"
$ws.Range("H9").Value = "bullet_damage"
$ws.Range("I9").Value = "bullet_velocity"
$ws.Range("J9").Value = "fire_rate"
$ws.Range("K9").Value = "price"
$ws.Range("L9").Value = "0st"
$ws.Range("M9").Value = "100st"
$ws.Range("N9").Value = "200st"
$ws.Range("O9").Value = "avg"
$ws.Range("P9").Value = "vel loss"
$ws.Range("Q9").Value = "suppression"
$ws.Range("R9").Value = "pen"
$ws.Range("S9").Value = "strength"
$ws.Range("A10").Value = "5.56x45_federal_american_eagle_training_xm193blx_55gr_fmj"
$ws.Range("B10").Value = "5.56x45 Federal American Eagle Training 55gr FMJ (XM193BLX)"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0.1
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0.04
$ws.Range("I10").Value = 150
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 55
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 44
$ws.Range("O10").Formula = "=AVERAGE(L10:N10)"
$ws.Range("P10").Value = 40
$ws.Range("Q10").Value = 30
$ws.Range("R10").Value = 0.8
$ws.Range("S10").Formula = "=C10-D10*20-E10*0.8-F10-0.6-G10*5+I10/200+(O10-50)"
$ws.Range("A11").Value = "5.56x45_cbcdefense_m196_55gr_red_tracer"
$ws.Range("B11").Value = "5.56x45 CBC Defense M196 55gr Red Tracer"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0.09
$ws.Range("E11").Value = -1
$ws.Range("F11").Value = -2
$ws.Range("G11").Value = 0.1
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 175
$ws.Range("K11").Value = 600
$ws.Range("L11").Value = 53
$ws.Range("M11").Value = 48
$ws.Range("N11").Value = 42
$ws.Range("O11").Formula = "=AVERAGE(L11:N11)"
$ws.Range("P11").Value = 40
$ws.Range("Q11").Value = 60
$ws.Range("R11").Value = 0.6
$ws.Range("S11").Formula = "=C11-D11*20-E11*0.8-F11-0.6-G11*5+I11/200+(O11-50)"
$ws.Range("A12").Value = "5.56x45_winchester_m855_62gr_fmj_greentip_vm855"
$ws.Range("B12").Value = "5.56x45 Winchester M855 62gr FMJ Green Tip (VM855)"
$ws.Range("C12").Value = -2
$ws.Range("D12").Value = 0.15
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = -0.05
$ws.Range("H12").Value = -0.02
$ws.Range("I12").Value = 100
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 58
$ws.Range("M12").Value = 55
$ws.Range("N12").Value = 51
$ws.Range("O12").Formula = "=AVERAGE(L12:N12)"
$ws.Range("P12").Value = 30
$ws.Range("Q12").Value = 50
$ws.Range("R12").Value = 0.7
$ws.Range("S12").Formula = "=C12-D12*20-E12*0.8-F12-0.6-G12*5+I12/200+(O12-50)"
$ws.Range("A13").Value = "5.56x45_cbcdefense_sat_ip_62gr_clf_aep_97"
$ws.Range("B13").Value = "5.56x45 CBC Defense SAT IP 62gr CLF (AEP-97)"
$ws.Range("C13").Value = -1
$ws.Range("D13").Value = 0.14000000000000001
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = -0.2
$ws.Range("H13").Value = 0.1
$ws.Range("I13").Value = 50
$ws.Range("K13").Value = 2000
$ws.Range("L13").Value = 58
$ws.Range("M13").Value = 53
$ws.Range("N13").Value = 46
$ws.Range("O13").Formula = "=AVERAGE(L13:N13)"
$ws.Range("P13").Value = 20
$ws.Range("Q13").Value = 50
$ws.Range("R13").Value = 1
$ws.Range("S13").Formula = "=C13-D13*20-E13*0.8-F13-0.6-G13*5+I13/200+(O13-50)"
$ws.Range("A14").Value = "5.56x45_novx_copper_pentagon_55gr_mchp_556n55cp"
$ws.Range("B14").Value = "5.56x45 NOVX Copper Pentagon 55gr MCHP (556N55CP)"
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 0.08
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 3
$ws.Range("G14").Value = 0.3
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 200
$ws.Range("K14").Value = 1200
$ws.Range("L14").Value = 65
$ws.Range("M14").Value = 53
$ws.Range("N14").Value = 40
$ws.Range("O14").Formula = "=AVERAGE(L14:N14)"
$ws.Range("P14").Value = 80
$ws.Range("Q14").Value = 40
$ws.Range("R14").Value = 0.3
$ws.Range("S14").Formula = "=C14-D14*20-E14*0.8-F14-0.6-G14*5+I14/200+(O14-50)"

# --- View / selection settings ---
$ws.Range("N21").Select()
$win = $excel.ActiveWindow
$win.Zoom = 130
